$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5683
$ws1.Range("F5").Value = 958
$ws1.Range("F6").Value = 36

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5683
$ws4.Range("F5").Value = 958
$ws4.Range("F6").Value = 36
